$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct/clean up the species labels in the "Morro Bay" block (rows 3-17).
$ws.Range("B6").Value  = "Jack mackerel "
$ws.Range("B10").Value = "Perch"
$ws.Range("B11").Value = "White seabass"
$ws.Range("B3").Value  = "Albacore"
$ws.Range("B5").Value  = "Rockfish"
$ws.Range("B7").Value  = "Giant Pacific oyster"
$ws.Range("B8").Value  = "Petrale sole"
$ws.Range("B9").Value  = "Salmon"
$ws.Range("B13").Value = "English sole"
$ws.Range("B14").Value = "California halibut"

# Leave the final selection on B15, matching where editing ended.
$ws.Range("B15").Select()

$wb.Save()
